$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fix the typo "Coordinador:COORDINADOR ZONA 4.4" (missing space after the
# colon) -> "Coordinador: COORDINADOR ZONA 4.4" on every cell that carries
# that text. Excel/WPS de-duplicates the shared-string table on save, so
# this also folds these cells onto the (already correctly spaced) string
# used by E37, letting the old duplicate entry disappear.
$fixedText = "Coordinador: COORDINADOR ZONA 4.4"
$cellsToFix = @("E2","E3","E4","E6","E11","E12","E13","E14","E15","E17","E19","E21","E22","E24","E26","E27","E28","E30","E35")
foreach ($addr in $cellsToFix) {
    $ws.Range($addr).Value = $fixedText
}

# Reflect the author's final cursor/viewport position captured in the
# saved workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E42").Select()
